$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 (Dom Pérignon) — rename product, fix country, bump modified date
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "Dom Pérignon Vintage Champagne3"
$ws.Range("D2").Value = "Dom Pérignon Vintage Champagne3"
$ws.Range("N2").Value = "France"
$ws.Range("AU2").Value = "2025-03-29T07:26:26.162Z"

# ---------------------------------------------------------------------------
# Row 3 (Macallan Rare Cask) — rename product, normalize volume, bump date
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = "Macallan Rare Cask Single Malt2"
$ws.Range("D3").Value = "Macallan Rare Cask Single Malt2"
$ws.Range("AC3").Value = "750ML"
$ws.Range("AU3").Value = "2025-03-29T07:36:19.940Z"

# ---------------------------------------------------------------------------
# Row 16 (Casamigos Blanco Tequila) — rename product, convert the boolean
# flag columns + alcohol% to plain text/number, normalize volume, add image
# + modified-date columns
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = "Casamigos Blanco Tequila3"
$ws.Range("D16").Value = "Casamigos Blanco Tequila3"

function Set-TextValue($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue "W16" "true"
Set-TextValue "X16" "false"
Set-TextValue "Y16" "false"
Set-TextValue "Z16" "true"

$ws.Range("AA16").Value = 40

$ws.Range("AC16").Value = "750ML"

Set-TextValue "AS16" "true"

$ws.Range("AT16").Value = "https://ext.same-assets.com/1701767421/1355704146.jpeg"
$ws.Range("AU16").Value = "2025-03-29T08:15:38.675Z"

# ---------------------------------------------------------------------------
# Remove the two trailing test rows (608 / 609) — shrinks the used range
# from A1:AX27 down to A1:AX25
# ---------------------------------------------------------------------------
$ws.Range("A26:A27").EntireRow.Delete()
